$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before the old column G (old G:J shift to I:L)
$ws.Columns("G:H").Insert()

# New column G: header + formula D/(1000*64000) for rows 2-25
$ws.Range("G1").Value = "mol/dL"
$ws.Range("G2").Formula = "=D2/(1000*64000)"
$ws.Range("G3:G25").Formula = "=D3/(1000*64000)"

# Autofit the new column to its content, matching the bestFit column width behavior
$ws.Columns("G:G").AutoFit()

# Restore the selection/active cell as left by the author
$ws.Range("H15").Select()
